$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: the politeness_score was stored as text "3"; it should be the
# real number 3.
$ws.Range("B34").Value = 3

# Append a new annotation row (row 35).
$ws.Range("A35").Value = "Sunsi Wu"

# B35's politeness_score must stay textual "3" (as authored), unlike B34
# above. Setting .Value = "3" directly gets auto-coerced to the number 3
# by Excel, so instead enter it with a leading apostrophe (forces text),
# then round-trip it through a values-only copy/paste so the cell keeps
# its text type without carrying over the apostrophe's quote-prefix
# formatting. The scratch cell is cleared afterwards.
$ws.Range("Z1").Value = "'3"
$ws.Range("Z1").Copy()
$ws.Range("B35").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("C35").Value = "do not"
$ws.Range("D35").Value = "DFT"
$ws.Range("E35").Value = "MET"
$ws.Range("F35").Value = "b3917550-3902-443d-ae6f-4c206bcc883a"
$ws.Range("G35").Value = "HkJ1rgbCb_annotated.xlsx"
$ws.Range("H35").Value = "However, these selections do not seem to directly incorporate the competing/augmenting effects of having different subgraphs within a molecule."
